$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03835166666666667
$ws.Range("H2").Value = 0.115055
$ws.Range("I2").Value = 0.0181239951898282
$ws.Range("J2").Value = 0.0181239951898282
$ws.Range("M2").Value = 0.7893693333333335
$ws.Range("N2").Value = 2.368108
$ws.Range("O2").Value = 0.1840020898203156
$ws.Range("P2").Value = 0.1840020898203156
$ws.Range("Q2").Value = 0.03027362954888889
$ws.Range("R2").Value = 0.2724626659400001
$ws.Range("S2").Value = 0.003334852990821735
$ws.Range("T2").Value = 0.003334852990821735
# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03835166666666667
$ws.Range("H3").Value = 0.115055
$ws.Range("I3").Value = 0.0181239951898282
$ws.Range("J3").Value = 0.0181239951898282
$ws.Range("O3").Value = 0.6904048063380857
$ws.Range("P3").Value = 0.6904048063380857
$ws.Range("Q3").Value = 0.1135914236966667
$ws.Range("R3").Value = 1.02232281327
$ws.Range("S3").Value = 0.01251289338910573
$ws.Range("T3").Value = 0.01251289338910573
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03835166666666667
$ws.Range("H4").Value = 0.115055
$ws.Range("I4").Value = 0.0181239951898282
$ws.Range("J4").Value = 0.0181239951898282
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5387946666666666
$ws.Range("N4").Value = 1.616384
$ws.Range("O4").Value = 0.1255931038415988
$ws.Range("P4").Value = 0.1255931038415988
$ws.Range("Q4").Value = 0.02066367345777778
$ws.Range("R4").Value = 0.18597306112
$ws.Range("S4").Value = 0.00227624880990073
$ws.Range("T4").Value = 0.00227624880990073
# Row 5
$ws.Range("I5").Value = 0.3727881574250648
$ws.Range("J5").Value = 0.3727881574250648
$ws.Range("M5").Value = 0.7893693333333335
$ws.Range("N5").Value = 2.368108
$ws.Range("O5").Value = 0.1840020898203156
$ws.Range("P5").Value = 0.1840020898203156
$ws.Range("Q5").Value = 0.6226911042457779
$ws.Range("R5").Value = 5.604219938212001
$ws.Range("S5").Value = 0.06859380002647671
$ws.Range("T5").Value = 0.06859380002647671
# Row 6
$ws.Range("I6").Value = 0.3727881574250648
$ws.Range("J6").Value = 0.3727881574250648
$ws.Range("O6").Value = 0.6904048063380857
$ws.Range("P6").Value = 0.6904048063380857
$ws.Range("S6").Value = 0.2573747356321836
$ws.Range("T6").Value = 0.2573747356321836
# Row 7
$ws.Range("I7").Value = 0.3727881574250648
$ws.Range("J7").Value = 0.3727881574250648
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5387946666666666
$ws.Range("N7").Value = 1.616384
$ws.Range("O7").Value = 0.1255931038415988
$ws.Range("P7").Value = 0.1255931038415988
$ws.Range("Q7").Value = 0.4250261972195555
$ws.Range("R7").Value = 3.825235774976
$ws.Range("S7").Value = 0.04681962176640445
$ws.Range("T7").Value = 0.04681962176640445
# Row 8
$ws.Range("G8").Value = 1.288873333333333
$ws.Range("H8").Value = 3.86662
$ws.Range("I8").Value = 0.6090878473851071
$ws.Range("J8").Value = 0.609087847385107
$ws.Range("M8").Value = 0.7893693333333335
$ws.Range("N8").Value = 2.368108
$ws.Range("O8").Value = 0.1840020898203156
$ws.Range("P8").Value = 0.1840020898203156
$ws.Range("Q8").Value = 1.017397083884445
$ws.Range("R8").Value = 9.15657375496
$ws.Range("S8").Value = 0.1120734368030172
$ws.Range("T8").Value = 0.1120734368030171
# Row 9
$ws.Range("G9").Value = 1.288873333333333
$ws.Range("H9").Value = 3.86662
$ws.Range("I9").Value = 0.6090878473851071
$ws.Range("J9").Value = 0.609087847385107
$ws.Range("O9").Value = 0.6904048063380857
$ws.Range("P9").Value = 0.6904048063380857
$ws.Range("Q9").Value = 3.817434015853333
$ws.Range("R9").Value = 34.35690614268
$ws.Range("S9").Value = 0.4205171773167963
$ws.Range("T9").Value = 0.4205171773167963
# Row 10
$ws.Range("G10").Value = 1.288873333333333
$ws.Range("H10").Value = 3.86662
$ws.Range("I10").Value = 0.6090878473851071
$ws.Range("J10").Value = 0.609087847385107
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5387946666666666
$ws.Range("N10").Value = 1.616384
$ws.Range("O10").Value = 0.1255931038415988
$ws.Range("P10").Value = 0.1255931038415988
$ws.Range("Q10").Value = 0.6944380780088889
$ws.Range("R10").Value = 6.249942702079998
$ws.Range("S10").Value = 0.07649723326529366
$ws.Range("T10").Value = 0.07649723326529365

Write-Host "Applied TPM updates"